# Applies the cryptos.xlsx update described in the commit
# "Updated cryptos list on Sun Sep 15 03:40:10 UTC 2024 with GitHub Actions"
# All edited cells hold plain text (inline/shared strings) in the source workbook,
# so for values that look numeric we briefly force a Text number format before
# assigning them (otherwise Excel would silently coerce them to numbers and drop
# things like trailing zeros), then restore the original (Normal) style so the
# cell formatting matches the original file again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.165.81'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '2.424.40'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.25%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.80%  '
$ws.Range("E9").Value = '  -0.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.148'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.355'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("D14").Value = '2.857.21'
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("D15").Value = '60.097.36'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("D17").Value = '2.401.87'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '329.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("E24").Value = '  +3.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.59%  '
$ws.Range("D28").Value = '0.0₃0776'
$ws.Range("E28").Value = '  -1.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.78'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.60'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.81%  '
$ws.Range("E33").Value = '  -3.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.33'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.23'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '332.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.41%  '
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.84'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '146.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("E44").Value = '  +2.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0968'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("E47").Value = '  +0.85%  '
$ws.Range("E48").Value = '  -1.23%  '
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("E50").Value = '  -2.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.66'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.07%  '
